$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 = Idaho. Fill in the new results from this run.
$ws.Range("B36").Value = 44022
$ws.Range("B36").NumberFormat = "YYYY-MM-DD"

$ws.Range("C36").Value = 9928
$ws.Range("D36").Value = 101
$ws.Range("E36").Value = 145
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 1.46
$ws.Range("H36").Value = 0.99

$ws.Range("J36").Value = $true

$ws.Range("O36").Value = "Success!"
